$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$newRows = @(
    @(45791.341157407405, 11, 6, 366, 545, 518, 575, 4305, 575, 2841, 277, 519, 30, 4919, 6551),
    @(45791.348599537036, 11, 6, 366, 545, 518, 575, 4305, 575, 2841, 277, 519, 30, 4919, 6551),
    @(45792.00105324074,  11, 6, 375, 558, 531, 575, 4305, 575, 2841, 277, 519, 30, 5032, 6630),
    @(45795.334745370368, 11, 6, 379, 558, 531, 620, 4751, 620, 3881, 369, 580, 30, 5221, 6695),
    @(45795.86577546296,  11, 6, 379, 558, 531, 671, 5349, 671, 672,  49,  655, 30, 5466, 6701),
    @(45799.25136574074,  11, 6, 383, 579, 552, 673, 5361, 673, 50,   4,   656, 30, 5530, 6775),
    @(45799.260034722225, 11, 6, 383, 579, 552, 673, 5361, 673, 50,   4,   656, 30, 5530, 6775),
    @(45802.923344907409, 11, 6, 385, 606, 578, 682, 5472, 682, 237,  12,  673, 30, 5635, 6845)
)

$startRow = 66
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
